$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column H header ---
$ws.Range("H1").Value = "p_adj"

# --- Fill new column H (p_adj) values ---
$ws.Range("H2").Value = 0.00026478379040417699
$ws.Range("H3").Value = 0.037592766515825202
$ws.Range("H4").Value = 0.00000025556498250978198
$ws.Range("H5").Value = 0.43725182114501898
$ws.Range("H6").Value = 0.77943452842727501
$ws.Range("H7").Value = 0.43725182114501898
$ws.Range("H8").Value = 0.0015550608938807999
$ws.Range("H9").Value = 0.00026478379040417699
$ws.Range("H10").Value = 0.43725182114501898
$ws.Range("H11").Value = 0.0080647194116229195
$ws.Range("H12").Value = 0.00043295081772978198
$ws.Range("H13").Value = 0.000101512145040505

# --- Move the "significant" highlight (yellow fill) from column E to column H ---
# First, copy the highlighted style (taken from a currently-highlighted E cell) onto the
# H rows that should now be highlighted.
$highlightRows = @(2,3,4,8,9,11,12,13)
$ws.Range("E2").Copy()
foreach ($r in $highlightRows) {
  $ws.Range("H$r").PasteSpecial(-4122)
}

# Then clear the old highlight from column E (copy the plain style from D2, which has no fill).
$ws.Range("D2").Copy()
foreach ($r in $highlightRows) {
  $ws.Range("E$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# --- Column A width (best-fit-like width) ---
$ws.Columns("A").ColumnWidth = 17.5

# --- View: zoom + selection on the new column ---
$excel.ActiveWindow.Zoom = 115
$ws.Range("H11:H13").Select()
